$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.562.55"
$ws.Range("E2").Value = "  -2.33%  "
$ws.Range("D3").Value = "1.583.23"
$ws.Range("E3").Value = "  -2.96%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.72%  "
$ws.Range("E6").Value = "  -2.15%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.248"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.41%  "
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.50"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0831"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("D12").Value = "1.805.00"
$ws.Range("E12").Value = "  -2.94%  "
$ws.Range("D13").Value = "1.576.38"
$ws.Range("E13").Value = "  -3.43%  "
$ws.Range("E14").Value = "  -1.80%  "
$ws.Range("E15").Value = "  -2.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.31"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").Value = "26.579.43"
$ws.Range("E17").Value = "  -2.08%  "
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "208.53"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.71"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.98%  "
$ws.Range("E22").Value = "  -3.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.38"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.41%  "
$ws.Range("E24").Value = "  -2.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.43"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("E26").Value = "  +1.90%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  -4.72%  "
$ws.Range("E29").Value = "  -1.97%  "
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("E32").Value = "  -3.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.670"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +24.09%  "
$ws.Range("E34").Value = "  -2.89%  "
$ws.Range("D35").Value = "1.320.13"
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("E36").Value = "  -4.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.42"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.00%  "
$ws.Range("E38").Value = "  -1.04%  "
$ws.Range("E39").Value = "  -3.47%  "
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("E41").Value = "  -2.70%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.17"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.01%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.29"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.99"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("D45").Value = "1.718.51"
$ws.Range("E45").Value = "  -2.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.95"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("E47").Value = "  +0.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.830"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0986"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0507"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.49"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.31%  "
